$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1) - column F ("想去人数") updates
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F6").Value  = 37
$wsExhibit.Range("F7").Value  = 590
$wsExhibit.Range("F8").Value  = 110
$wsExhibit.Range("F9").Value  = 8729
$wsExhibit.Range("F11").Value = 330
$wsExhibit.Range("F12").Value = 1144
$wsExhibit.Range("F13").Value = 983
$wsExhibit.Range("F16").Value = 5
$wsExhibit.Range("F18").Value = 258
$wsExhibit.Range("F21").Value = 1028

# Sheet "全部类型" (sheet4) - column F ("想去人数") updates
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F7").Value  = 37
$wsAll.Range("F9").Value  = 590
$wsAll.Range("F10").Value = 110
$wsAll.Range("F11").Value = 8729
$wsAll.Range("F13").Value = 330
$wsAll.Range("F14").Value = 1144
$wsAll.Range("F15").Value = 983
$wsAll.Range("F18").Value = 5
$wsAll.Range("F20").Value = 258
$wsAll.Range("F23").Value = 1028
